$d = $word.ActiveDocument

$d.Content.Find.Execute("40×34=", $true, $false, $false, $false, $false, $true, 1, $false, "47×14=", 2) | Out-Null
$d.Content.Find.Execute("38×49=", $true, $false, $false, $false, $false, $true, 1, $false, "42×27=", 2) | Out-Null
$d.Content.Find.Execute("80×86=", $true, $false, $false, $false, $false, $true, 1, $false, "45×96=", 2) | Out-Null
$d.Content.Find.Execute("52×27=", $true, $false, $false, $false, $false, $true, 1, $false, "38×81=", 2) | Out-Null
$d.Content.Find.Execute("67×75=", $true, $false, $false, $false, $false, $true, 1, $false, "68×72=", 2) | Out-Null
$d.Content.Find.Execute("98×63=", $true, $false, $false, $false, $false, $true, 1, $false, "91×11=", 2) | Out-Null
$d.Content.Find.Execute("56×30=", $true, $false, $false, $false, $false, $true, 1, $false, "72×46=", 2) | Out-Null
$d.Content.Find.Execute("21×90=", $true, $false, $false, $false, $false, $true, 1, $false, "42×70=", 2) | Out-Null
$d.Content.Find.Execute("23×88=", $true, $false, $false, $false, $false, $true, 1, $false, "94×76=", 2) | Out-Null
$d.Content.Find.Execute("58×77=", $true, $false, $false, $false, $false, $true, 1, $false, "85×40=", 2) | Out-Null
$d.Content.Find.Execute("68×82=", $true, $false, $false, $false, $false, $true, 1, $false, "90×53=", 2) | Out-Null
$d.Content.Find.Execute("58×47=", $true, $false, $false, $false, $false, $true, 1, $false, "59×64=", 2) | Out-Null
$d.Content.Find.Execute("24×52=", $true, $false, $false, $false, $false, $true, 1, $false, "24×27=", 2) | Out-Null
$d.Content.Find.Execute("23×55=", $true, $false, $false, $false, $false, $true, 1, $false, "31×84=", 2) | Out-Null
$d.Content.Find.Execute("63×23=", $true, $false, $false, $false, $false, $true, 1, $false, "80×23=", 2) | Out-Null
$d.Content.Find.Execute("34×75=", $true, $false, $false, $false, $false, $true, 1, $false, "97×94=", 2) | Out-Null
$d.Content.Find.Execute("89×34=", $true, $false, $false, $false, $false, $true, 1, $false, "99×93=", 2) | Out-Null
$d.Content.Find.Execute("67×78=", $true, $false, $false, $false, $false, $true, 1, $false, "29×90=", 2) | Out-Null
$d.Content.Find.Execute("37×92=", $true, $false, $false, $false, $false, $true, 1, $false, "91×66=", 2) | Out-Null
$d.Content.Find.Execute("53×65=", $true, $false, $false, $false, $false, $true, 1, $false, "82×72=", 2) | Out-Null
$d.Content.Find.Execute("62×33=", $true, $false, $false, $false, $false, $true, 1, $false, "40×57=", 2) | Out-Null
$d.Content.Find.Execute("28×45=", $true, $false, $false, $false, $false, $true, 1, $false, "44×79=", 2) | Out-Null
$d.Content.Find.Execute("14×52=", $true, $false, $false, $false, $false, $true, 1, $false, "67×89=", 2) | Out-Null
$d.Content.Find.Execute("33×59=", $true, $false, $false, $false, $false, $true, 1, $false, "57×55=", 2) | Out-Null
$d.Content.Find.Execute("69×33=", $true, $false, $false, $false, $false, $true, 1, $false, "94×59=", 2) | Out-Null
